$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by
# one day (45838 -> 45839, i.e. 2025-06-30 -> 2025-07-01) for every data
# row (rows 2 through 43).
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45838) {
        $cell.Value2 = 45839
    }
}
